# Swap the colour scheme carried by the deck's theme (ppt/theme/theme1.xml,
# the theme wired to the slide master / all regular slides) from the
# "Integral" palette over to the stock "Office Theme" palette, and push the
# "Integral" palette onto the notes-master theme (ppt/theme/theme2.xml) so
# the two themes trade places.

$p = $ppt.ActivePresentation

# ---- Target (stock "Office Theme") colour values, as 0x00BBGGRR ints ----
# dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
# accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
# accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

# ---- Original ("Integral") colour values, as 0x00BBGGRR ints ----
# dk1=000000 lt1=FFFFFF dk2=455F51 lt2=E3DED1
# accent1=99CB38 accent2=63A537 accent3=E6D024 accent4=CC9700
# accent5=4EB3CF accent6=378DA6 hlink=6B9F25 folHlink=B26B02
$integralColors = @(0, 16777215, 5332805, 13754083, 3722137, 3646819, 2412774, 38860, 13611854, 10915127, 2465643, 158642)

# The slide master's theme (ppt/theme/theme1.xml) drives every normal
# slide's look -- recolour it from Integral to the Office Theme palette.
$masterTheme = $p.SlideMaster.Theme
$masterTheme.Name = "Office Theme"
$masterScheme = $masterTheme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Item($i).RGB = $officeColors[$i - 1]
}

# The notes master's theme (ppt/theme/theme2.xml) takes on the palette the
# slide master used to carry (Integral), completing the swap.
$notesTheme = $p.NotesMaster.Theme
$notesTheme.Name = "Integral"
$notesScheme = $notesTheme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Item($i).RGB = $integralColors[$i - 1]
}
